$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.347.72'
$ws.Range('E2').Value = '  -2.28%  '
$ws.Range('D3').Value = '3.489.62'
$ws.Range('E3').Value = '  -2.30%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').Value = "'611.52"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +4.70%  '
$ws.Range('D6').Value = "'185.76"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.37%  '
$ws.Range('E7').Value = '  +0.93%  '
$ws.Range('D8').Value = "'0.999"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.08%  '
$ws.Range('D9').Value = "'0.213"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.26%  '
$ws.Range('D10').Value = "'0.651"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.39%  '
$ws.Range('D11').Value = "'53.18"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.83%  '
$ws.Range('E12').Value = '  -4.28%  '
$ws.Range('D13').Value = "'9.60"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.03%  '
$ws.Range('D14').Value = '4.049.36'
$ws.Range('E14').Value = '  -1.97%  '
$ws.Range('D15').Value = "'606.56"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +8.02%  '
$ws.Range('D16').Value = '69.413.03'
$ws.Range('E16').Value = '  -2.14%  '
$ws.Range('D17').Value = "'12.65"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.71%  '
$ws.Range('D18').Value = "'18.86"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.46%  '
$ws.Range('D19').Value = '3.504.42'
$ws.Range('E19').Value = '  -0.60%  '
$ws.Range('E20').Value = '  -0.44%  '
$ws.Range('D21').Value = "'0.985"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.17%  '
$ws.Range('D23').Value = "'104.69"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E24').Value = '  +1.33%  '
$ws.Range('D25').Value = "'5.01"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.80%  '
$ws.Range('D26').Value = "'3.03"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.44%  '
$ws.Range('D27').Value = "'10.96"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.19%  '
$ws.Range('D28').Value = "'9.93"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +8.31%  '
$ws.Range('D29').Value = "'33.69"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E30').Value = '  -3.95%  '
$ws.Range('D31').Value = "'12.49"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.20%  '
$ws.Range('D32').Value = "'0.117"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.06%  '
$ws.Range('D33').Value = "'3.90"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +15.91%  '
$ws.Range('E34').Value = '  -0.65%  '
$ws.Range('D35').Value = "'3.17"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -6.71%  '
$ws.Range('D36').Value = "'0.998"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.19%  '
$ws.Range('D37').Value = "'523.51"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -5.65%  '
$ws.Range('E38').Value = '  -5.63%  '
$ws.Range('D39').Value = '3.570.53'
$ws.Range('E39').Value = '  +0.28%  '
$ws.Range('D40').Value = "'3.59"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +4.15%  '
$ws.Range('D41').Value = "'36.70"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.14%  '
$ws.Range('D42').Value = '0.0₃0774'
$ws.Range('E42').Value = '  -3.72%  '
$ws.Range('E43').Value = '  +0.96%  '
$ws.Range('E44').Value = '  +2.72%  '
$ws.Range('D45').Value = "'2.98"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.26%  '
$ws.Range('E46').Value = '  +5.60%  '
$ws.Range('D47').Value = "'3.33"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -5.06%  '
$ws.Range('D48').Value = "'8.86"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -5.53%  '
$ws.Range('E49').Value = '  +0.32%  '
$ws.Range('D50').Value = "'131.23"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.98%  '
$ws.Range('D51').Value = "'1.36"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -9.52%  '
